# This script re-shuffles the per-fund data rows across the nine
# "compliance by fund" tables in this workbook so that each table's
# row 2,3,4,6,7,8 carries a different fund's figures (row 5 / SDVD is
# unaffected). Only the literal (non-formula) cell contents move; any
# formula cell is left completely untouched so that it keeps
# referencing its own row and recalculates naturally against the
# newly-placed values.
#
# The row permutation (new row <- old row) is:
#   2 <- 6   (TDVI's numbers move up into row 2)
#   3 <- 4   (FGSI's numbers move up into row 3)
#   4 <- 7   (RDVI's numbers move up into row 4)
#   5 <- 5   (unchanged)
#   6 <- 2   (KNG's numbers move down into row 6)
#   7 <- 8   (DOGG's numbers move up into row 7)
#   8 <- 3   (FDND's numbers move down into row 8)

$wb = $excel.ActiveWorkbook

# Sheet name -> last data column (as a column index, A=1)
$sheetCols = @{
    "Prospectus_80pct"                 = 16  # A..P
    "40Act_Diversification"            = 23  # A..W
    "IRS_Diversification"              = 18  # A..R
    "Illiquid"                         = 8   # A..H
    "Real_Estate"                      = 4   # A..D
    "Commodities"                      = 3   # A..C
    "12d1_Other_Investment_Companies"  = 10  # A..J
    "12d2_Insurance_Companies"         = 5   # A..E
    "12d3_Securities_Business"         = 11  # A..K
}

# new row -> old row (where its replacement content comes from)
$rowMap = @{2=6; 3=4; 4=7; 5=5; 6=2; 7=8; 8=3}

foreach ($sheetName in $sheetCols.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $maxCol = $sheetCols[$sheetName]

    # Snapshot every source row's non-formula cell values (column B
    # onward - column A is just the report date and is identical on
    # every row already). Must snapshot everything up front since the
    # permutation includes a 4-cycle (3<-4<-7<-8<-3), not just simple
    # pairwise swaps.
    $snapshot = @{}
    foreach ($r in @(2,3,4,5,6,7,8)) {
        $rowVals = @{}
        for ($c = 2; $c -le $maxCol; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if (-not $cell.HasFormula) {
                $rowVals[$c] = $cell.Value2
            }
        }
        $snapshot[$r] = $rowVals
    }

    # Write the snapshot back out according to the permutation.
    foreach ($destRow in @(2,3,4,5,6,7,8)) {
        $srcRow = $rowMap[$destRow]
        $rowVals = $snapshot[$srcRow]
        foreach ($c in $rowVals.Keys) {
            $ws.Cells.Item($destRow, $c).Value2 = $rowVals[$c]
        }
    }
}
